# Insert a new data row at row 465 (pushes existing rows 465:485 down to 466:486)
# and populate it with the new weekly price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(465).Insert()

$ws.Cells.Item(465, 1).Value  = 5
$ws.Cells.Item(465, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(465, 3).Value  = "Maule"
$ws.Cells.Item(465, 4).Value  = 44939
$ws.Cells.Item(465, 5).Value  = 7
$ws.Cells.Item(465, 6).Value  = 100112032
$ws.Cells.Item(465, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(465, 8).Value  = "Sin especificar"
$ws.Cells.Item(465, 9).Value  = "Primera"
$ws.Cells.Item(465, 10).Value = 500
$ws.Cells.Item(465, 11).Value = 7000
$ws.Cells.Item(465, 12).Value = 7000
$ws.Cells.Item(465, 13).Value = 7000
$ws.Cells.Item(465, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(465, 15).Value = "Región del Maule"
$ws.Cells.Item(465, 16).Value = 140
$ws.Cells.Item(465, 17).Value = 50
$ws.Cells.Item(465, 18).Value = "Hortaliza"
